# SCD0011-029 - Penyelia SRM Mengakses Menu Report - Menu Product Holding Ratio - Report
# "Update Excel SCD0011 until SCD0016"
#
# 1. Rename worksheet SCD0198 -> SCD0011
# 2. Update TC_ID cell (B2) DGS-213 -> SCD0011-029
# 3. Widen column B to fit the new (longer) TC_ID text
# 4. Move the active selection from C2 to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "SCD0011"

# 2. Update the TC_ID value in B2
$ws.Cells.Item(2, 2).Value = "SCD0011-029"

# 3. Widen column B (was 9 characters, now fits "SCD0011-029")
$ws.Columns.Item(2).ColumnWidth = 11.7

# 4. Update the selected cell shown in the sheet view
$ws.Range("B3").Select()
